$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = 42587.819618055553

$ws.Range("B5").Value = "Noun"
$ws.Range("C5").Value = 13362
$ws.Range("D5").Value = 9857
$ws.Range("E5").Value = 1762
$ws.Range("F5").Value = 233
$ws.Range("G5").Value = 99
$ws.Range("H5").Value = 69
$ws.Range("I5").Value = 29
$ws.Range("J5").Value = 7
$ws.Range("K5").Value = 4
$ws.Range("L5").Value = 63
$ws.Range("M5").Value = 36
